$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008 and 2009 rows (rows 2 and 3); this shifts all subsequent
# rows (2010..2020 data) up by two, turning old rows 4..14 into new rows 2..12.
$ws.Rows("2:3").Delete()

# Append the new 2021 row (new row 13) of data, copying the year-label cell
# formatting (bold / bordered / centered style) from the row above it first.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 33069.410479219
$ws.Range("C13").Value = 26309.163394618
$ws.Range("D13").Value = 59378.573873837
$ws.Range("E13").Value = 350551.725728905
$ws.Range("F13").Value = 222500.404022872
$ws.Range("G13").Value = 381539.241347715
$ws.Range("H13").Value = 125163.693454687
$ws.Range("I13").Value = 217160.753925738
$ws.Range("J13").Value = 14995.7754615991
$ws.Range("K13").Value = 73885.8578782973
$ws.Range("L13").Value = 81208.2162489895
$ws.Range("M13").Value = 21.1749863022
$ws.Range("N13").Value = 118576.166042561
$ws.Range("O13").Value = 311829.891861099
$ws.Range("P13").Value = 112651.881122121
$ws.Range("Q13").Value = 703434.34609804
$ws.Range("R13").Value = 246908.232733487
$ws.Range("S13").Value = 253973.53064218
$ws.Range("T13").Value = 249744.128219165
$ws.Range("U13").Value = 412863.014835705
$ws.Range("V13").Value = 1333788.51218235
$ws.Range("W13").Value = 2116898.59145877
$ws.Range("X13").Value = 55509.0806292732
$ws.Range("Y13").Value = 5339.6500971341
$ws.Range("Z13").Value = 3430395.12805284
$ws.Range("AA13").Value = 3430395.12805284
$ws.Range("AB13").Value = 1032441.18381299
$ws.Range("AC13").Value = 412951.554706329
$ws.Range("AD13").Value = 556618.198436346
$ws.Range("AE13").Value = 2002010.93695566

